# Updates cryptos list prices / volume(1h) figures, and swaps two pairs of
# rows (Toncoin/InjectiveProtocol and Monero/ARBITRUM) back to their
# original relative order with refreshed figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (single "." as decimal
# point) must be forced to Text format first, otherwise Excel's COM layer
# would coerce the string into a numeric value and silently drop
# formatting such as trailing zeros (e.g. "72.80" -> 72.8).

$ws.Range("D2").Value = '42.784.44'
$ws.Range("E2").Value = '  -7.84%  '
$ws.Range("D3").Value = '2.515.11'
$ws.Range("E3").Value = '  -4.12%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '296.28'
$ws.Range("E5").Value = '  -4.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.36'
$ws.Range("E6").Value = '  -6.61%  '
$ws.Range("E7").Value = '  -5.41%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.546'
$ws.Range("E9").Value = '  -6.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.06'
$ws.Range("E10").Value = '  -8.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0802'
$ws.Range("E11").Value = '  -5.27%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.64'
$ws.Range("E13").Value = '  +0.55%  '
$ws.Range("D14").Value = '2.901.51'
$ws.Range("E14").Value = '  -3.84%  '
$ws.Range("D15").Value = '2.520.38'
$ws.Range("E15").Value = '  -3.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.864'
$ws.Range("E16").Value = '  -6.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.08'
$ws.Range("E17").Value = '  -6.29%  '
$ws.Range("D18").Value = '42.780.96'
$ws.Range("E18").Value = '  -8.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.56'
$ws.Range("E19").Value = '  -3.65%  '
$ws.Range("E20").Value = '  -5.42%  '
$ws.Range("E21").Value = '  -6.09%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.80'
$ws.Range("E22").Value = '  +1.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '258.98'
$ws.Range("E23").Value = '  -6.37%  '
$ws.Range("E24").Value = '  -5.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.18'
$ws.Range("E25").Value = '  -1.34%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '28.89'
$ws.Range("E26").Value = '  -2.50%  '
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.94'
$ws.Range("E28").Value = '  -6.80%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.12'
$ws.Range("E29").Value = '  -7.08%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.64'
$ws.Range("E30").Value = '  -5.35%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.96'
$ws.Range("E31").Value = '  -7.71%  '
$ws.Range("E32").Value = '  -4.78%  '
$ws.Range("B33").Value = 'Monero'
$ws.Range("C33").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '151.68'
$ws.Range("E33").Value = '  -0.32%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.19'
$ws.Range("E34").Value = '  -3.43%  '
$ws.Range("E35").Value = '  -2.92%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0798'
$ws.Range("E36").Value = '  -5.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.115'
$ws.Range("E37").Value = '  -6.30%  '
$ws.Range("E38").Value = '  -3.95%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '23.65'
$ws.Range("E39").Value = '  -2.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.31'
$ws.Range("E40").Value = '  +1.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.45'
$ws.Range("E41").Value = '  -4.61%  '
$ws.Range("E42").Value = '  -7.05%  '
$ws.Range("E43").Value = '  -6.30%  '
$ws.Range("D44").Value = '2.021.90'
$ws.Range("E44").Value = '  -5.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '85.58'
$ws.Range("E46").Value = '  -9.97%  '
$ws.Range("E47").Value = '  +3.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.88'
$ws.Range("E48").Value = '  -6.73%  '
$ws.Range("D49").Value = '2.761.87'
$ws.Range("E49").Value = '  -3.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '102.82'
$ws.Range("E50").Value = '  -6.39%  '
$ws.Range("E51").Value = '  -7.83%  '
